$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'327.77"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'-0.25%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'44.52"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'0.79%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'5.402"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'-1.91%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.08322"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'3.01%"
$ws.Range("E5").Style = "Normal"
$ws.Range("B6").Value = "FTXToken"
$ws.Range("C6").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D6").Value = "'1.926"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'-4.84%"
$ws.Range("E6").Style = "Normal"
$ws.Range("B7").Value = "MXToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D7").Value = "'0.9729"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'2.23%"
$ws.Range("E7").Style = "Normal"
$ws.Range("B8").Value = "BTSEToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D8").Value = "'2.527"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'-3.58%"
$ws.Range("E8").Style = "Normal"
$ws.Range("B9").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C9").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D9").Value = "'0.1112"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'-1.75%"
$ws.Range("E9").Style = "Normal"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").Value = "'0.1885"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'0.64%"
$ws.Range("E10").Style = "Normal"
$ws.Range("B11").Value = "MCDex"
$ws.Range("C11").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D11").Value = "'9.196"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'-8.77%"
$ws.Range("E11").Style = "Normal"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").Value = "'0.09687"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'-2.88%"
$ws.Range("E12").Style = "Normal"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").Value = "'0.04698"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'-2.89%"
$ws.Range("E13").Style = "Normal"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").Value = "'0.1065"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'0.71%"
$ws.Range("E14").Style = "Normal"
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").Value = "'0.001284"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'2.42%"
$ws.Range("E15").Style = "Normal"
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").Value = "'0.006118"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'1.91%"
$ws.Range("E16").Style = "Normal"
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").Value = "'3.385"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'0.39%"
$ws.Range("E17").Style = "Normal"
$ws.Range("B18").Value = "GateToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D18").Value = "'4.427"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'0.38%"
$ws.Range("E18").Style = "Normal"
$ws.Range("B19").Value = "BitpandaEcosystemToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D19").Value = "'0.3323"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'0.87%"
$ws.Range("E19").Style = "Normal"
$ws.Range("B20").Value = "ProBitToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("D20").Value = "'0.1372"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'-1.94%"
$ws.Range("E20").Style = "Normal"
$ws.Range("B21").Value = "ZBToken"
$ws.Range("C21").Value = "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
$ws.Range("D21").Value = "'0.2729"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'6.13%"
$ws.Range("E21").Style = "Normal"
$ws.Range("B22").Value = "CoinExToken"
$ws.Range("C22").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D22").Value = "'0.04178"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'2.39%"
$ws.Range("E22").Style = "Normal"
$ws.Range("B23").Value = "BitKan"
$ws.Range("C23").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("D23").Value = "'0.001300"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'-0.51%"
$ws.Range("E23").Style = "Normal"
$ws.Range("B24").Value = "HotbitToken"
$ws.Range("C24").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("D24").Value = "'0.004425"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'1.39%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.0001307"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'4.58%"
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'-20.10%"
$ws.Range("E26").Style = "Normal"
$ws.Range("D38").Value = "'0.02635"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'1.60%"
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'0.05594"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'-0.66%"
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.008011"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'5.09%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.1414"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'1.04%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.007413"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'0.96%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.002128"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'7.59%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.008658"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'5.11%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.3364"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Value = "'0.00006858"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'-3.19%"
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "'0.59%"
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").Value = "'0.32%"
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'0.003525"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'0.35%"
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.003549"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'1.50%"
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.00002112"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'0.59%"
$ws.Range("E51").Style = "Normal"
